$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2:H14 values from 50 to 60
$ws.Range("H2:H14").Value = 60

# Update the selection shown in the sheet view
$ws.Range("H2:H14").Select()
